$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the new row by shifting existing rows down, then rewrite
# every row in final row-major order so the shared-string table is
# rebuilt in the same order the cells are encountered.
$ws.Rows.Item(1).Insert()

$ws.Range("A1").Value = "Copper Horse - Warhorse Blend (12oz)"
$ws.Range("B1").Value = 7

$ws.Range("A2").Value = "Copper Horse - Rumble Pony (12oz)"
$ws.Range("B2").Value = 9

$ws.Range("A3").Value = "Copper Horse - Clocktower Espresso (12oz)"
$ws.Range("B3").Value = 14

$ws.Range("A4").Value = "Copper Horse - Carriage House Blend (12oz)"
$ws.Range("B4").Value = 9

$ws.Range("A5").Value = "Copper Horse - Sleigh Bells (12oz)"
$ws.Range("B5").Value = 22

$ws.Range("A6").Value = "Copper Horse - Warhorse 5lb"
$ws.Range("B6").Value = 6
